$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the two leave dates (2025-03-25 and 2025-03-27 -> rows 26 and 28) as "Planned Leave"
$ws.Range("B26").Value = "Planned Leave"
$ws.Range("B28").Value = "Planned Leave"

# Append the Leave Summary block below the existing data (rows 35-39)
$ws.Range("A35").Value = "Mar-2025 Leave Summary"

$ws.Range("A36").Value = "Sick Leave"
$ws.Range("B36").Value = 0

$ws.Range("A37").Value = "Planned Leave"
$ws.Range("B37").Value = 3

$ws.Range("A38").Value = "Unplanned Leave"
$ws.Range("B38").Value = 0

$ws.Range("A39").Value = "WFH"
$ws.Range("B39").Value = 0

# Column A grew wider to fit the new summary labels
$ws.Columns.Item(1).AutoFit()
